$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fijación de precios" (4° semestre, rows 3-8): every UNIDAD row had its own
# distinct CLAVE (MCZ27402..MCZ27406). Excel's correct key is the same one
# used for the subject's first unit (row 3) - fix rows 4-8 to match A3.
$ws.Range("A4").Value = "MCZ27401 "
$ws.Range("A5").Value = "MCZ27401 "
$ws.Range("A6").Value = "MCZ27401 "
$ws.Range("A7").Value = "MCZ27401 "
$ws.Range("A8").Value = "MCZ27401 "

# "Segmentación de mercados" (4° semestre, rows 9-14): same fix, using A9's
# key (MCI27401) for rows 10-14.
$ws.Range("A10").Value = "MCI27401"
$ws.Range("A11").Value = "MCI27401"
$ws.Range("A12").Value = "MCI27401"
$ws.Range("A13").Value = "MCI27401"
$ws.Range("A14").Value = "MCI27401"

# "Psicología del consumidor" (4° semestre, rows 15-20): same fix, using
# A15's key (MIM27401) for rows 16-20.
$ws.Range("A16").Value = "MIM27401"
$ws.Range("A17").Value = "MIM27401"
$ws.Range("A18").Value = "MIM27401"
$ws.Range("A19").Value = "MIM27401"
$ws.Range("A20").Value = "MIM27401"

# Leave the selection where the author ended up after making the fix.
$ws.Range("A16:A20").Select()
